# Applies the edit described by the diff:
#  - Adds 3 new worksheets: bt_3, cg_3, n_conv3 (each with the same
#    Row/FD1/FD2 x h=1e-2..h=1e-12/Exact layout used by the existing sheets)
#  - Updates the existing "time_3" sheet's column widths and data values
#
# Helper: the engine snaps ColumnWidth (characters) onto a 1/6-character
# grid (stored = (round(input*6)+5)/6). We invert that so the written
# "characters" width lands as close as possible to the widths recorded in
# the target workbook.
function Set-ClosestColumnWidth($ws, $colIndex, $targetStoredWidth) {
    $px = [double]$targetStoredWidth * 6.0 - 5.0
    $pxFloor = [Math]::Floor($px)
    $pxCeil = [Math]::Ceiling($px)
    $inputFloor = $pxFloor / 6.0
    $inputCeil = $pxCeil / 6.0
    $storedFloor = ([Math]::Floor($inputFloor * 6.0 + 0.5) + 5.0) / 6.0
    $storedCeil = ([Math]::Floor($inputCeil * 6.0 + 0.5) + 5.0) / 6.0
    if ([Math]::Abs($storedFloor - $targetStoredWidth) -le [Math]::Abs($storedCeil - $targetStoredWidth)) {
        $ws.Columns.Item($colIndex).ColumnWidth = $inputFloor
    } else {
        $ws.Columns.Item($colIndex).ColumnWidth = $inputCeil
    }
}

function Fill-ResultSheet($ws, $row1Label, $row2Label, $headers, $row2Vals, $row3Vals, $colWidths) {
    # Header row
    $ws.Cells.Item(1,1).Value = $row1Label
    for ($c = 0; $c -lt $headers.Length; $c++) {
        $ws.Cells.Item(1, $c + 2).Value = $headers[$c]
    }
    # Data rows
    $ws.Cells.Item(2,1).Value = "FD1"
    for ($c = 0; $c -lt $row2Vals.Length; $c++) {
        $ws.Cells.Item(2, $c + 2).Value = $row2Vals[$c]
    }
    $ws.Cells.Item(3,1).Value = "FD2"
    for ($c = 0; $c -lt $row3Vals.Length; $c++) {
        $ws.Cells.Item(3, $c + 2).Value = $row3Vals[$c]
    }
    # Column widths
    for ($c = 0; $c -lt $colWidths.Length; $c++) {
        Set-ClosestColumnWidth $ws ($c + 1) $colWidths[$c]
    }
}

$wb = $excel.ActiveWorkbook

$headers = @("h=1e-2","h=1e-4","h=1e-6","h=1e-8","h=1e-10","h=1e-12","Exact")

# ---------------------------------------------------------------------
# Update existing sheet "time_3" (B2:H3 values + column widths C/D)
# ---------------------------------------------------------------------
$wsTime3 = $wb.Worksheets.Item("time_3")

$time3Row2 = @(0.009599509999999999,0.00924808,0.0095053875,0.009072255555555554,0.008918344444444445,0.009417433333333334,0.021389933333333336)
$time3Row3 = @(0.044539528571428566,0.043413455555555557,0.04436099,0.044145809999999994,0.04322984444444446,0.04229206666666666,0.021389933333333336)

for ($c = 0; $c -lt $time3Row2.Length; $c++) {
    $wsTime3.Cells.Item(2, $c + 2).Value = $time3Row2[$c]
}
for ($c = 0; $c -lt $time3Row3.Length; $c++) {
    $wsTime3.Cells.Item(3, $c + 2).Value = $time3Row3[$c]
}

Set-ClosestColumnWidth $wsTime3 3 10.7109375
Set-ClosestColumnWidth $wsTime3 4 12.7109375

# ---------------------------------------------------------------------
# Add new sheet "bt_3" after the last existing sheet (ExactComparison)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsBt3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsBt3.Name = "bt_3"

$bt3Row2 = @(0.5015730807035156,0.4314632468980295,0.3685216052063879,0.46114859375728934,0.46114859375728934,0.46114859375728934,0.46114859375728934)
$bt3Row3 = @(0.42518552875695736,0.41514874014874015,0.41624585559368177,0.45023559773559774,0.46114859375728934,0.46114859375728934,0.46114859375728934)
$bt3ColWidths = @(4.82421875,12.7109375,12.7109375,12.7109375,12.7109375,12.7109375,12.7109375,12.7109375)

Fill-ResultSheet $wsBt3 "Row" "FD1" $headers $bt3Row2 $bt3Row3 $bt3ColWidths

# ---------------------------------------------------------------------
# Add new sheet "cg_3" after "bt_3"
# ---------------------------------------------------------------------
$wsCg3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsBt3)
$wsCg3.Name = "cg_3"

$cg3Row2 = @(13.21628845389715,13.134415186263013,13.245764018590107,13.56649996541301,13.56649996541301,13.56649996541301,13.56649996541301)
$cg3Row3 = @(13.02899196042053,12.853978428978429,13.220391094412836,13.479692807192809,13.56649996541301,13.56649996541301,13.56649996541301)
$cg3ColWidths = @(4.82421875,11.7109375,11.7109375,11.7109375,11.7109375,11.7109375,11.7109375,11.7109375)

Fill-ResultSheet $wsCg3 "Row" "FD1" $headers $cg3Row2 $cg3Row3 $cg3ColWidths

# ---------------------------------------------------------------------
# Add new sheet "n_conv3" after "cg_3"
# ---------------------------------------------------------------------
$wsNConv3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCg3)
$wsNConv3.Name = "n_conv3"

$nconv3Row2 = @(10,10,8,9,9,9,9)
$nconv3Row3 = @(7,9,10,10,9,9,9)
$nconv3ColWidths = @(4.82421875,6.82421875,6.82421875,6.82421875,6.82421875,7.82421875,7.82421875,5.6015625)

Fill-ResultSheet $wsNConv3 "Row" "FD1" $headers $nconv3Row2 $nconv3Row3 $nconv3ColWidths
